# Update the cryptocurrency price/volume table (columns D and E)
# to reflect the latest scraped values, matching the commit
# "Updated cryptos list on Sun Feb 26 23:37:27 UTC 2023 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell without Excel silently
# reinterpreting number-looking strings (e.g. "6.590", "0.9993") as
# real numbers, which would drop formatting such as trailing zeros.
# We temporarily force a Text number format, assign the value, then
# restore the cell to the default "Normal" style so no visible
# formatting is left behind.
function Set-TextValue {
    param($cell, [string]$text, [bool]$forceText)
    if ($forceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

# Row 2: D2: '23.501.17' -> '23.540.93'; E2: '  +1.72%  ' -> '  +1.89%  '
Set-TextValue $ws.Cells.Item(2, 4) '23.540.93' $false
Set-TextValue $ws.Cells.Item(2, 5) '  +1.89%  ' $false

# Row 3: D3: '1.639.98' -> '1.640.46'; E3: '  +3.06%  ' -> '  +3.10%  '
Set-TextValue $ws.Cells.Item(3, 4) '1.640.46' $false
Set-TextValue $ws.Cells.Item(3, 5) '  +3.10%  ' $false

# Row 4: D4: '1.001' -> '0.9993'; E4: '  +0.37%  ' -> '  +0.16%  '
Set-TextValue $ws.Cells.Item(4, 4) '0.9993' $true
Set-TextValue $ws.Cells.Item(4, 5) '  +0.16%  ' $false

# Row 5: D5: '308.44' -> '308.78'; E5: '  +2.28%  ' -> '  +2.45%  '
Set-TextValue $ws.Cells.Item(5, 4) '308.78' $true
Set-TextValue $ws.Cells.Item(5, 5) '  +2.45%  ' $false

# Row 6: D6: '0.9996' -> '0.9998'; E6: '  +0.24%  ' -> '  +0.27%  '
Set-TextValue $ws.Cells.Item(6, 4) '0.9998' $true
Set-TextValue $ws.Cells.Item(6, 5) '  +0.27%  ' $false

# Row 7: D7: '0.3770' -> '0.3775'; E7: '  +0.26%  ' -> '  +0.35%  '
Set-TextValue $ws.Cells.Item(7, 4) '0.3775' $true
Set-TextValue $ws.Cells.Item(7, 5) '  +0.35%  ' $false

# Row 8: D8: '52.86' -> '52.84'; E8: '  +3.45%  ' -> '  +3.56%  '
Set-TextValue $ws.Cells.Item(8, 4) '52.84' $true
Set-TextValue $ws.Cells.Item(8, 5) '  +3.56%  ' $false

# Row 9: D9: '0.3683' -> '0.3685'; E9: '  +2.05%  ' -> '  +2.15%  '
Set-TextValue $ws.Cells.Item(9, 4) '0.3685' $true
Set-TextValue $ws.Cells.Item(9, 5) '  +2.15%  ' $false

# Row 10: D10: '1.277' -> '1.279'; E10: '  +2.54%  ' -> '  +2.59%  '
Set-TextValue $ws.Cells.Item(10, 4) '1.279' $true
Set-TextValue $ws.Cells.Item(10, 5) '  +2.59%  ' $false

# Row 11: D11: '0.08211' -> '0.08217'; E11: '  +2.13%  ' -> '  +2.20%  '
Set-TextValue $ws.Cells.Item(11, 4) '0.08217' $true
Set-TextValue $ws.Cells.Item(11, 5) '  +2.20%  ' $false

# Row 12: D12: '1.001' -> '0.9994'; E12: '  +0.40%  ' -> '  +0.19%  '
Set-TextValue $ws.Cells.Item(12, 4) '0.9994' $true
Set-TextValue $ws.Cells.Item(12, 5) '  +0.19%  ' $false

# Row 13: D13: '23.21' -> '23.26'; E13: '  +3.54%  ' -> '  +3.91%  '
Set-TextValue $ws.Cells.Item(13, 4) '23.26' $true
Set-TextValue $ws.Cells.Item(13, 5) '  +3.91%  ' $false

# Row 14: D14: '6.668' -> '6.674'; E14: '  +2.27%  ' -> '  +2.35%  '
Set-TextValue $ws.Cells.Item(14, 4) '6.674' $true
Set-TextValue $ws.Cells.Item(14, 5) '  +2.35%  ' $false

# Row 15: D15: '0.00001281' -> '0.00001285'; E15: '  +3.52%  ' -> '  +3.86%  '
Set-TextValue $ws.Cells.Item(15, 4) '0.00001285' $true
Set-TextValue $ws.Cells.Item(15, 5) '  +3.86%  ' $false

# Row 16: D16: '7.469' -> '7.478'; E16: '  +1.44%  ' -> '  +1.58%  '
Set-TextValue $ws.Cells.Item(16, 4) '7.478' $true
Set-TextValue $ws.Cells.Item(16, 5) '  +1.58%  ' $false

# Row 17: D17: '1.639.94' -> '1.640.56'; E17: '  +2.67%  ' -> '  +3.10%  '
Set-TextValue $ws.Cells.Item(17, 4) '1.640.56' $false
Set-TextValue $ws.Cells.Item(17, 5) '  +3.10%  ' $false

# Row 18: E18: '  +2.32%  ' -> '  +2.38%  '
Set-TextValue $ws.Cells.Item(18, 5) '  +2.38%  ' $false

# Row 19: E19: '  +2.90%  ' -> '  +2.86%  '
Set-TextValue $ws.Cells.Item(19, 5) '  +2.86%  ' $false

# Row 20: D20: '18.41' -> '18.42'; E20: '  +2.73%  ' -> '  +2.77%  '
Set-TextValue $ws.Cells.Item(20, 4) '18.42' $true
Set-TextValue $ws.Cells.Item(20, 5) '  +2.77%  ' $false

# Row 21: D21: '6.585' -> '6.590'
Set-TextValue $ws.Cells.Item(21, 4) '6.590' $true

# Row 22: D22: '0.9968' -> '0.9981'; E22: '  -0.06%  ' -> '  +0.07%  '
Set-TextValue $ws.Cells.Item(22, 4) '0.9981' $true
Set-TextValue $ws.Cells.Item(22, 5) '  +0.07%  ' $false

# Row 23: D23: '23.495.46' -> '23.548.40'; E23: '  +1.74%  ' -> '  +1.90%  '
Set-TextValue $ws.Cells.Item(23, 4) '23.548.40' $false
Set-TextValue $ws.Cells.Item(23, 5) '  +1.90%  ' $false

# Row 24: D24: '12.96' -> '12.97'; E24: '  +1.19%  ' -> '  +1.40%  '
Set-TextValue $ws.Cells.Item(24, 4) '12.97' $true
Set-TextValue $ws.Cells.Item(24, 5) '  +1.40%  ' $false

# Row 25: D25: '3.116' -> '3.100'; E25: '  +6.59%  ' -> '  +6.43%  '
Set-TextValue $ws.Cells.Item(25, 4) '3.100' $true
Set-TextValue $ws.Cells.Item(25, 5) '  +6.43%  ' $false

# Row 26: D26: '2.410' -> '2.412'; E26: '  +1.14%  ' -> '  +1.34%  '
Set-TextValue $ws.Cells.Item(26, 4) '2.412' $true
Set-TextValue $ws.Cells.Item(26, 5) '  +1.34%  ' $false

# Row 27: D27: '21.43' -> '21.41'; E27: '  +2.51%  ' -> '  +2.38%  '
Set-TextValue $ws.Cells.Item(27, 4) '21.41' $true
Set-TextValue $ws.Cells.Item(27, 5) '  +2.38%  ' $false

# Row 28: D28: '151.56' -> '151.69'; E28: '  +2.06%  ' -> '  +2.05%  '
Set-TextValue $ws.Cells.Item(28, 4) '151.69' $true
Set-TextValue $ws.Cells.Item(28, 5) '  +2.05%  ' $false

# Row 29: D29: '5.334' -> '5.324'; E29: '  +2.85%  ' -> '  +2.69%  '
Set-TextValue $ws.Cells.Item(29, 4) '5.324' $true
Set-TextValue $ws.Cells.Item(29, 5) '  +2.69%  ' $false

# Row 30: E30: '  +2.38%  ' -> '  +2.40%  '
Set-TextValue $ws.Cells.Item(30, 5) '  +2.40%  ' $false

# Row 31: D31: '2.421' -> '2.420'; E31: '  +2.07%  ' -> '  +2.06%  '
Set-TextValue $ws.Cells.Item(31, 4) '2.420' $true
Set-TextValue $ws.Cells.Item(31, 5) '  +2.06%  ' $false

# Row 32: D32: '6.866' -> '6.862'; E32: '  +2.31%  ' -> '  +2.18%  '
Set-TextValue $ws.Cells.Item(32, 4) '6.862' $true
Set-TextValue $ws.Cells.Item(32, 5) '  +2.18%  ' $false

# Row 33: D33: '1.820.58' -> '1.821.16'; E33: '  +3.06%  ' -> '  +2.82%  '
Set-TextValue $ws.Cells.Item(33, 4) '1.821.16' $false
Set-TextValue $ws.Cells.Item(33, 5) '  +2.82%  ' $false

# Row 34: D34: '0.9776' -> '0.9815'; E34: '  +2.13%  ' -> '  +2.91%  '
Set-TextValue $ws.Cells.Item(34, 4) '0.9815' $true
Set-TextValue $ws.Cells.Item(34, 5) '  +2.91%  ' $false

# Row 35: D35: '0.02815' -> '0.02820'; E35: '  +5.33%  ' -> '  +5.37%  '
Set-TextValue $ws.Cells.Item(35, 4) '0.02820' $true
Set-TextValue $ws.Cells.Item(35, 5) '  +5.37%  ' $false

# Row 36: E36: '  +3.58%  ' -> '  +3.43%  '
Set-TextValue $ws.Cells.Item(36, 5) '  +3.43%  ' $false

# Row 37: D37: '0.07476' -> '0.07492'; E37: '  +0.16%  ' -> '  +0.30%  '
Set-TextValue $ws.Cells.Item(37, 4) '0.07492' $true
Set-TextValue $ws.Cells.Item(37, 5) '  +0.30%  ' $false

# Row 38: D38: '6.222' -> '6.227'; E38: '  +2.11%  ' -> '  +2.22%  '
Set-TextValue $ws.Cells.Item(38, 4) '6.227' $true
Set-TextValue $ws.Cells.Item(38, 5) '  +2.22%  ' $false

# Row 39: D39: '0.2542' -> '0.2549'; E39: '  +1.74%  ' -> '  +1.99%  '
Set-TextValue $ws.Cells.Item(39, 4) '0.2549' $true
Set-TextValue $ws.Cells.Item(39, 5) '  +1.99%  ' $false

# Row 40: E40: '  +1.13%  ' -> '  +1.12%  '
Set-TextValue $ws.Cells.Item(40, 5) '  +1.12%  ' $false

# Row 41: D41: '1.403' -> '1.401'; E41: '  +3.41%  ' -> '  +3.27%  '
Set-TextValue $ws.Cells.Item(41, 4) '1.401' $true
Set-TextValue $ws.Cells.Item(41, 5) '  +3.27%  ' $false

# Row 42: D42: '0.7164' -> '0.7176'; E42: '  +0.99%  ' -> '  +1.24%  '
Set-TextValue $ws.Cells.Item(42, 4) '0.7176' $true
Set-TextValue $ws.Cells.Item(42, 5) '  +1.24%  ' $false

# Row 43: D43: '12.65' -> '12.67'; E43: '  +3.86%  ' -> '  +3.92%  '
Set-TextValue $ws.Cells.Item(43, 4) '12.67' $true
Set-TextValue $ws.Cells.Item(43, 5) '  +3.92%  ' $false

# Row 44: D44: '16.17' -> '16.29'; E44: '  +8.16%  ' -> '  +8.14%  '
Set-TextValue $ws.Cells.Item(44, 4) '16.29' $true
Set-TextValue $ws.Cells.Item(44, 5) '  +8.14%  ' $false

# Row 45: D45: '0.6621' -> '0.6634'; E45: '  +2.03%  ' -> '  +2.02%  '
Set-TextValue $ws.Cells.Item(45, 4) '0.6634' $true
Set-TextValue $ws.Cells.Item(45, 5) '  +2.02%  ' $false

# Row 46: D46: '2.362' -> '2.371'; E46: '  +3.62%  ' -> '  +3.96%  '
Set-TextValue $ws.Cells.Item(46, 4) '2.371' $true
Set-TextValue $ws.Cells.Item(46, 5) '  +3.96%  ' $false

# Row 47: D47: '4.046' -> '4.049'; E47: '  +1.37%  ' -> '  +1.42%  '
Set-TextValue $ws.Cells.Item(47, 4) '4.049' $true
Set-TextValue $ws.Cells.Item(47, 5) '  +1.42%  ' $false

# Row 48: D48: '0.08068' -> '0.08074'; E48: '  +2.08%  ' -> '  +2.09%  '
Set-TextValue $ws.Cells.Item(48, 4) '0.08074' $true
Set-TextValue $ws.Cells.Item(48, 5) '  +2.09%  ' $false

# Row 49: D49: '130.90' -> '131.11'; E49: '  -0.41%  ' -> '  -0.10%  '
Set-TextValue $ws.Cells.Item(49, 4) '131.11' $true
Set-TextValue $ws.Cells.Item(49, 5) '  -0.10%  ' $false

# Row 50: D50: '1.221' -> '1.223'; E50: '  +0.19%  ' -> '  +0.59%  '
Set-TextValue $ws.Cells.Item(50, 4) '1.223' $true
Set-TextValue $ws.Cells.Item(50, 5) '  +0.59%  ' $false

# Row 51: D51: '1.218' -> '1.220'; E51: '  +2.68%  ' -> '  +2.75%  '
Set-TextValue $ws.Cells.Item(51, 4) '1.220' $true
Set-TextValue $ws.Cells.Item(51, 5) '  +2.75%  ' $false

Write-Output "Updated cryptos list on Sun Feb 26 23:37:27 UTC 2023 with GitHub Actions"
